$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.316.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "'2.998.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'564.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").Value = "'139.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.46%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("D9").Value = "'2.990.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("E10").Value = "  +10.18%  "
$ws.Range("D11").Value = "'4.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.38%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.01%  "
$ws.Range("E13").Value = "  +10.90%  "
$ws.Range("D14").Value = "'33.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.18%  "
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "'3.492.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.46%  "
$ws.Range("D18").Value = "'2.997.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").Value = "'59.236.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "'429.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.27%  "
$ws.Range("D21").Value = "'13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.68%  "
$ws.Range("D22").Value = "'0.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.46%  "
$ws.Range("D23").Value = "'7.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.10%  "
$ws.Range("D24").Value = "'13.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("D25").Value = "'80.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.11%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  +12.39%  "
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("D30").Value = "'7.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.32%  "
$ws.Range("D31").Value = "'25.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.95%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").Value = "'0.0994"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("E34").Value = "  +12.19%  "
$ws.Range("D35").Value = "'0.0₃0778"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +27.01%  "
$ws.Range("D36").Value = "'5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.50%  "
$ws.Range("E37").Value = "  +5.56%  "
$ws.Range("D38").Value = "'49.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("D39").Value = "'8.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("E40").Value = "  +17.44%  "
$ws.Range("D41").Value = "'406.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.38%  "
$ws.Range("D42").Value = "'0.0352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").Value = "'2.763.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.85%  "
$ws.Range("D44").Value = "'0.108"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").Value = "'0.248"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.36%  "
$ws.Range("D47").Value = "'124.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("D48").Value = "'2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.37%  "
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "'32.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +23.33%  "
$ws.Range("D51").Value = "'23.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.41%  "
